$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header from "Cruise" to "Month"
$ws.Range("A1").Value = "Month"

# Delete columns E:H (In_situ_DOU_mean, In_situ_DOU_sd, OPD_mean, OPD_sd)
$ws.Range("E1:H12").Delete()

# Replace cruise identifiers with month names
$ws.Range("A2:A6").Value = "March"
$ws.Range("A7:A12").Value = "October"
